$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (the most recent "Day" entry, 2020-06-04) has a "blog" widget in
# column I that pointed at blog post #151. Replace it with a reference to
# the newly published blog post #152 ("no more poverty").
$ws.Range("I11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 152"
